$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtered save games) for row 2
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 35.8657594772105
